{"js": "// Insert a new \"List Bullet\" paragraph right after the\n// \"Docente(s) Respons\u00e1vel(eis)\" heading, listing the three professors\n// (one per line, separated by manual line breaks / <w:br/>), matching\n// the target OOXML diff.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Find the \"Docente(s) Respons\u00e1vel(eis)\" heading paragraph.\nlet target = null;\nfor (const p of paragraphs.items) {\n  if (p.text.indexOf(\"Docente(s) Respons\u00e1vel(eis)\") !== -1) {\n    target = p;\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error('Could not find the \"Docente(s) Respons\u00e1vel(eis)\" paragraph');\n}\n\n// A zero-width insertion point right after the heading paragraph\n// (just past its paragraph mark, ahead of whatever follows it).\nconst insertionPoint = target.getRange(\"After\");\n\n// Build the new paragraph (style + three runs, the first two runs each\n// carrying a trailing manual line break) as a self-contained WordML\n// fragment wrapped in the Flat OPC envelope Word.Range.insertOoxml expects.\nconst flatOpcPackage =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" ' +\n  'pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  \"<w:body>\" +\n  \"<w:p>\" +\n  '<w:pPr><w:pStyle w:val=\"ListBullet\"/></w:pPr>' +\n  \"<w:r><w:t>7459752 - Maria Ismenia Sodero Toledo Faria</w:t><w:br/></w:r>\" +\n  \"<w:r><w:t>2166002 - Sandra Giacomin Schneider</w:t><w:br/></w:r>\" +\n  \"<w:r><w:t>1922320 - Sebastiao Ribeiro</w:t></w:r>\" +\n  \"</w:p>\" +\n  \"</w:body>\" +\n  \"</w:document>\" +\n  \"</pkg:xmlData>\" +\n  \"</pkg:part>\" +\n  \"</pkg:package>\";\n\ninsertionPoint.insertOoxml(flatOpcPackage, \"After\");\n\nawait context.sync();\n", "ps1": "# Insert a new \"List Bullet\" paragraph right after the\n# \"Docente(s) Respons\u00e1vel(eis)\" heading, listing the three professors\n# (one per line, separated by manual line breaks / <w:br/>), matching\n# the target OOXML diff.\n\n$d = $word.ActiveDocument\n\n# Find the \"Docente(s) Respons\u00e1vel(eis)\" heading paragraph.\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*Docente(s) Respons\u00e1vel(eis)*\") {\n        $target = $p\n        break\n    }\n}\n\nif ($null -eq $target) {\n    throw 'Could not find the \"Docente(s) Respons\u00e1vel(eis)\" paragraph'\n}\n\n# Split a brand-new (empty) paragraph in right after the heading, then\n# stamp the whole freshly-minted paragraph's range with the WordML\n# fragment holding the three names - the style plus the runs, the\n# first two runs each carrying a trailing manual line break.\n$target.Range.InsertParagraphAfter()\n$newPara = $target.Next()\n\n$fragment = '<w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:pPr><w:pStyle w:val=\"ListBullet\"/></w:pPr>' +\n    '<w:r><w:t>7459752 - Maria Ismenia Sodero Toledo Faria</w:t><w:br/></w:r>' +\n    '<w:r><w:t>2166002 - Sandra Giacomin Schneider</w:t><w:br/></w:r>' +\n    '<w:r><w:t>1922320 - Sebastiao Ribeiro</w:t></w:r>' +\n    '</w:p>'\n\n$newPara.Range.InsertXML($fragment)\n"}
